$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
$tfs = $theme.ThemeFontScheme
Write-Host "Name:" $tfs.Name
$major = $tfs.MajorFont
$minor = $tfs.MinorFont
Write-Host "Major Latin:" $major.Latin
Write-Host "Major EastAsian:" $major.EastAsian
Write-Host "Major ComplexScript:" $major.ComplexScript
Write-Host "Minor Latin:" $minor.Latin
Write-Host "Minor EastAsian:" $minor.EastAsian
Write-Host "Minor ComplexScript:" $minor.ComplexScript
